$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(2, "Morumbi - SP"),
    @(3, "Alphaville - SP"),
    @(4, "Home Office"),
    @(5, "Alphaville - SP"),
    @(6, "Home Office"),
    @(7, "Home Office"),
    @(8, "Home office"),
    @(9, "Home Office"),
    @(10, "home office"),
    @(11, "São Paulo"),
    @(12, "São Paulo"),
    @(13, "São Paulo"),
    @(14, "São Paulo"),
    @(15, "São Paulo"),
    @(16, "Home Office"),
    @(17, "Home office"),
    @(18, "Home Office"),
    @(19, "Home Office"),
    @(20, "Home Office"),
    @(21, "Home Office"),
    @(22, "Home Office"),
    @(23, "Home Office"),
    @(24, "Home Office"),
    @(25, "Home Office"),
    @(26, "Home Office"),
    @(27, "Home Office"),
    @(28, "Home Office"),
    @(29, "São Paulo"),
    @(30, "São Paulo"),
    @(31, "São Judas - SP"),
    @(32, "São Judas - SP"),
    @(33, "Home Office"),
    @(34, "Home Office"),
    @(35, "Home Office"),
    @(36, "Home Office"),
    @(37, "São Paulo")
)

foreach ($item in $values) {
    $row = $item[0]
    $val = $item[1]
    $ws.Cells.Item($row, 2).Value = $val
}

$ws.Cells.Item(2, 3).Value = "Teste"
